$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-09-08 -> 2023-09-09, i.e. serial 45177 -> 45178) for every data row
# (rows 2 through 57).
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 3).Value = 45178
}
